$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H98").Value = 6327.5
$ws.Range("I98").Value = 7272.6924
$ws.Range("J98").Value = 3255.625
$ws.Range("K98").Value = 7272.6924
$ws.Range("L98").Value = 3255.625
$ws.Range("M98").Value = -5774.6924
$ws.Range("N98").Value = -6251.625
$ws.Range("H122").Value = 6327.5
$ws.Range("I122").Value = 7272.6924
$ws.Range("J122").Value = 3255.625
$ws.Range("K122").Value = 21818.0772
$ws.Range("L122").Value = 9766.875
$ws.Range("M122").Value = -19368.0772
$ws.Range("N122").Value = -14666.875
$ws.Range("H132").Value = 2079.1482
$ws.Range("I132").Value = 1960.6842
$ws.Range("J132").Value = 2360.5
$ws.Range("K132").Value = 5882.0526
$ws.Range("L132").Value = 7081.5
$ws.Range("M132").Value = -3352.0526
$ws.Range("N132").Value = -12141.5
$ws.Range("H137").Value = 27736
$ws.Range("I137").Value = 27736
$ws.Range("K137").Value = 83208
$ws.Range("M137").Value = -80658
$ws.Range("M20").ClearContents()
$ws.Range("M35").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2210.3914
$ws.Range("I2").Value = 1090.4
$ws.Range("J2").Value = 3071.923
$ws.Range("K2").Value = 1090.4
$ws.Range("L2").Value = 3071.923
$ws.Range("M2").Value = -977.4000000000001
$ws.Range("N2").Value = -3297.923
$ws.Range("H32").Value = 2966.5
$ws.Range("I32").Value = 3018.7778
$ws.Range("J32").Value = 614
$ws.Range("K32").Value = 3018.7778
$ws.Range("L32").Value = 614
$ws.Range("M32").Value = -2731.7778
$ws.Range("N32").Value = -1188
$ws.Range("H61").Value = 3999.7144
$ws.Range("I61").Value = 3629.3333
$ws.Range("K61").Value = 3629.3333
$ws.Range("M61").Value = -3417.3333
$ws.Range("H74").Value = 2124.0833
$ws.Range("I74").Value = 2124.0833
$ws.Range("K74").Value = 2124.0833
$ws.Range("M74").Value = -1250.0833
$ws.Range("H77").Value = 2124.0833
$ws.Range("I77").Value = 2124.0833
$ws.Range("K77").Value = 10620.4165
$ws.Range("M77").Value = -6252.416499999999
$ws.Range("H116").Value = 2210.3914
$ws.Range("I116").Value = 1090.4
$ws.Range("J116").Value = 3071.923
$ws.Range("K116").Value = 1090.4
$ws.Range("L116").Value = 3071.923
$ws.Range("M116").Value = 1203.6
$ws.Range("N116").Value = -7659.923
$ws.Range("H136").Value = 3999.7144
$ws.Range("I136").Value = 3629.3333
$ws.Range("K136").Value = 10887.9999
$ws.Range("M136").Value = -8337.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2210.3914
$ws.Range("I3").Value = 1090.4
$ws.Range("J3").Value = 3071.923
$ws.Range("K3").Value = 1090.4
$ws.Range("L3").Value = 3071.923
$ws.Range("M3").Value = -976.4000000000001
$ws.Range("N3").Value = -3299.923
$ws.Range("H32").Value = 14750
$ws.Range("J32").Value = 14750
$ws.Range("L32").Value = 14750
$ws.Range("N32").Value = -15518
$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 4000
$ws.Range("K37").Value = 4000
$ws.Range("M37").Value = -3863
$ws.Range("H107").Value = 5040.88
$ws.Range("I107").Value = 1715.8572
$ws.Range("K107").Value = 1715.8572
$ws.Range("M107").Value = 204.1428000000001
$ws.Range("H134").Value = 3234.8462
$ws.Range("I134").Value = 3305.3
$ws.Range("K134").Value = 9915.900000000001
$ws.Range("M134").Value = -7380.900000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4629.5
$ws.Range("I31").Value = 1868.2632
$ws.Range("J31").Value = 15122.2
$ws.Range("K31").Value = 1868.2632
$ws.Range("L31").Value = 15122.2
$ws.Range("M31").Value = -1573.2632
$ws.Range("N31").Value = -15712.2
$ws.Range("H34").Value = 4629.5
$ws.Range("I34").Value = 1868.2632
$ws.Range("J34").Value = 15122.2
$ws.Range("K34").Value = 1868.2632
$ws.Range("L34").Value = 15122.2
$ws.Range("M34").Value = -1666.2632
$ws.Range("N34").Value = -15526.2
$ws.Range("H134").Value = 2629.2942
$ws.Range("I134").Value = 2752.6924
$ws.Range("J134").Value = 2228.25
$ws.Range("K134").Value = 8258.0772
$ws.Range("L134").Value = 6684.75
$ws.Range("M134").Value = -5723.0772
$ws.Range("N134").Value = -11754.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 47939.23
$ws.Range("I120").Value = 14610
$ws.Range("K120").Value = 43830
$ws.Range("M120").Value = -38992
$ws.Range("H121").Value = 733.26666
$ws.Range("I121").Value = 650.5
$ws.Range("K121").Value = 1951.5
$ws.Range("M121").Value = -641.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10075.236
$ws.Range("I70").Value = 7935.778
$ws.Range("K70").Value = 7935.778
$ws.Range("M70").Value = -7665.778
$ws.Range("H73").Value = 10075.236
$ws.Range("I73").Value = 7935.778
$ws.Range("K73").Value = 7935.778
$ws.Range("M73").Value = -6999.778
$ws.Range("H132").Value = 3510.389
$ws.Range("J132").Value = 4004.25
$ws.Range("L132").Value = 12012.75
$ws.Range("N132").Value = -17072.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 215997
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 215997
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 215997
$ws.Range("N129").Value = -225997
$ws.Range("M129").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2961571.5
$ws.Range("I2").Value = 4145000.2
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 4145000.2
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -4144888.2
$ws.Range("N2").Value = -3224
$ws.Range("H100").Value = 858.625
$ws.Range("I100").Value = 869.8333
$ws.Range("J100").Value = 825
$ws.Range("K100").Value = 1739.6666
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -1198.6666
$ws.Range("N100").Value = -2732
$ws.Range("H112").Value = 40599.8
$ws.Range("J112").Value = 40599.8
$ws.Range("L112").Value = 40599.8
$ws.Range("N112").Value = -43553.8
